$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for data rows 2 through 28
# from serial date 45501 (2024-07-28) to 45502 (2024-07-29).
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45501) {
        $cell.Value = 45502
    }
}
